$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Device")

# Remove the stray "Notes:" row (row 2) - it has no notes beneath it.
# Everything below shifts up by one row.
$ws.Rows("2").Delete()

# Update wording: "device(s)" -> "apparatus(es)"
$ws.Range("A1").Value = "This sheet summarizes the apparatuses connected to buses."
$ws.Range("B3").Value = "Type"
$ws.Range("C3").Value = "Parameters"

# Rename the sheet itself to match the new terminology.
$ws.Name = "Apparatus"

# Match the saved selection/cursor position.
[void]$ws.Range("C4").Select()
